# Update cryptocurrency price values in column D (Price) to reflect
# the latest scrape, per commit "Updated symbol list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'247.46"  # was 247.66
$ws.Range("D3").Value = "'22.40"  # was 22.37
$ws.Range("D4").Value = "'5.248"  # was 5.241
$ws.Range("D5").Value = "'0.05691"  # was 0.05688
$ws.Range("D8").Value = "'0.8079"  # was 0.8064
$ws.Range("D9").Value = "'0.8762"  # was 0.8813
$ws.Range("D11").Value = "'0.07441"  # was 0.07416
$ws.Range("D12").Value = "'0.03055"  # was 0.03056
$ws.Range("D14").Value = "'0.09395"  # was 0.09404
$ws.Range("D15").Value = "'3.878"  # was 3.883
$ws.Range("D17").Value = "'0.04769"  # was 0.04790
$ws.Range("D19").Value = "'0.006416"  # was 0.006417
$ws.Range("D20").Value = "'0.005038"  # was 0.005040
$ws.Range("D21").Value = "'0.0009960"  # was 0.0009966
$ws.Range("D22").Value = "'0.0001501"  # was 0.0001500
$ws.Range("D23").Value = "'3.692"  # was 3.690
$ws.Range("D24").Value = "'2.193"  # was 2.195
$ws.Range("D26").Value = "'0.1282"  # was 0.1352
$ws.Range("D27").Value = "'0.0004753"  # was 0.01827
$ws.Range("D40").Value = "'0.03952"  # was 0.03950
$ws.Range("D41").Value = "'0.006805"  # was 0.006810
$ws.Range("D42").Value = "'0.1064"  # was 0.1065
$ws.Range("D43").Value = "'0.003202"  # was 0.003200
$ws.Range("D44").Value = "'0.008451"  # was 0.008446
$ws.Range("D45").Value = "'0.00005585"  # was 0.00005587
$ws.Range("D48").Value = "'0.1553"  # was 0.2021
$ws.Range("D49").Value = "'0.00002101"  # was 0.00002100
$ws.Range("D50").Value = "'0.01011"  # was 0.01010
